# EV-103 Create config downloadable file for VT application, add party logo in the export
#
# Adds a new "image" column (F) to the PAR_2020_tab0a query table on Sheet1,
# containing the filename of each party's logo image, and expands the table /
# autofilter range to cover the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the query table by one column so the table range, autofilter and
# column metadata all expand to A1:F26.
$tbl = $ws.ListObjects.Item("PAR_2020_tab0a")
$newCol = $tbl.ListColumns.Add()

# Logo file names, in row order (rows 2-26 correspond to parties 1-25).
$logos = @(
    "slovenska-ludova-strana-andreja-hlinku.png",
    "dobra-volba.png",
    "sloboda-a-solidarita.png",
    "sme-rodina.png",
    "slovenske-hnutie-obrody.png",
    "za-ludi.png",
    "mame-toho-dost.png",
    "hlas-pravice.png",
    "slovenska-narodna-strana.png",
    "demokraticka-strana.png",
    "obycajni-ludia-a-nezavisle-osobnosti.png",
    "progresivne-slovensko-a-spolu.png",
    "starostovia-a-nezavisli-kandidati.png",
    "obciansky-hlas.png",
    "krestanskodemokraticke-hnutie.png",
    "slovenska-liga.png",
    "vlast.png",
    "most-hid.png",
    "smer-socialna-demokracia.png",
    "solidarita-hnutie-pracujucej-chudoby.png",
    "hlas-ludu.png",
    "madarska-komunitna-spolupatricnost.png",
    "praca-slovenskeho-naroda.png",
    "kotlebovci-ludova-strana-nase-slovensko.png",
    "socialisti.png"
)

for ($i = 0; $i -lt $logos.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $logos[$i]
}

# Header last, so the shared-string table order matches (data first, header last).
$ws.Range("F1").Value = "image"

# Match the saved window/selection state from the edit.
$ws.Range("G3").Select()
$excel.ActiveWindow.ScrollColumn = 3
